$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newDate = "2026-02-16 01:57:37"

# Remove existing hyperlinks up front; they will be rebuilt from scratch once
# every row is in its final place so the relationship ids come out in the
# same F2,F3,F4,F5 document order as a freshly written sheet.
$ws.Hyperlinks.Delete()

# Insert a new row above the current row 3 ("Salesforce" listing), shifting
# it (and the "Ticketmaster" listing below it) down by one row.
$ws.Rows.Item(3).Insert()

# Widen column D to fit the new, longer price range text. (31.17 is used
# instead of 32 to compensate for this runtime's char-width<->pixel
# round-trip rounding, so the persisted OOXML column width lands on an
# exact 32.)
$ws.Columns.Item(4).ColumnWidth = 31.17

# Row 2 keeps its original listing; only the scrape timestamp advances.
$ws.Range("A2").Value = $newDate

# Row 3: brand-new listing inserted by this scrape run.
$ws.Range("A3").Value = $newDate
$ws.Range("B3").Value = "地域情報サイト 店舗データ自動収集・一括管理システム構築"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5492383"
$ws.Range("G3").Value = 85
$ws.Range("H3").Value = "◇サイト"

# Row 4: previously row 3 ("Salesforce" listing); only the timestamp advances.
$ws.Range("A4").Value = $newDate

# Row 5: previously row 4 ("Ticketmaster" listing); only the timestamp advances.
$ws.Range("A5").Value = $newDate

# Rebuild the hyperlinks in top-to-bottom order so relationship ids land as
# rId1..rId4 matching F2..F5.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5251319")
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5492383")
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5492003")
$ws.Range("F4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5491983")
$ws.Range("F5").Style = "Hyperlink"
